# This script applies the diff: two new price-observation rows (Durazno,
# variety "Kurakata", date 44932) are inserted into the weekly dataset right
# before the current row 256, pushing the existing rows 256-309 down to
# 258-311. The sheet dimension grows from A1:T309 to A1:T311.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above current row 256 (shifts 256..309 -> 258..311)
$ws.Rows.Item(256).Insert()
$ws.Rows.Item(256).Insert()

# --- New row 256 ---
$ws.Cells.Item(256, 1).Value = 7
$ws.Cells.Item(256, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(256, 3).Value = "Ñuble"
$ws.Cells.Item(256, 4).Value = 44932
$ws.Cells.Item(256, 5).Value = 16
$ws.Cells.Item(256, 6).Value = "Fruta"
$ws.Cells.Item(256, 7).Value = 100103
$ws.Cells.Item(256, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(256, 9).Value = 100103004
$ws.Cells.Item(256, 10).Value = "Durazno"
$ws.Cells.Item(256, 11).Value = "Kurakata"
$ws.Cells.Item(256, 12).Value = "Especial"
$ws.Cells.Item(256, 13).Value = 80
$ws.Cells.Item(256, 14).Value = 17000
$ws.Cells.Item(256, 15).Value = 17000
$ws.Cells.Item(256, 16).Value = 17000
$ws.Cells.Item(256, 17).Value = "`$/caja 15 kilos granel"
$ws.Cells.Item(256, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(256, 19).Value = 1133
$ws.Cells.Item(256, 20).Value = 15

# --- New row 257 ---
$ws.Cells.Item(257, 1).Value = 7
$ws.Cells.Item(257, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(257, 3).Value = "Ñuble"
$ws.Cells.Item(257, 4).Value = 44932
$ws.Cells.Item(257, 5).Value = 16
$ws.Cells.Item(257, 6).Value = "Fruta"
$ws.Cells.Item(257, 7).Value = 100103
$ws.Cells.Item(257, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(257, 9).Value = 100103004
$ws.Cells.Item(257, 10).Value = "Durazno"
$ws.Cells.Item(257, 11).Value = "Kurakata"
$ws.Cells.Item(257, 12).Value = "Primera"
$ws.Cells.Item(257, 13).Value = 120
$ws.Cells.Item(257, 14).Value = 15000
$ws.Cells.Item(257, 15).Value = 16000
$ws.Cells.Item(257, 16).Value = 15500
$ws.Cells.Item(257, 17).Value = "`$/caja 15 kilos granel"
$ws.Cells.Item(257, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(257, 19).Value = 1033
$ws.Cells.Item(257, 20).Value = 15

# Make sure the two new date cells (column D) use the same date style/format
# as the rest of the column (numFmtId 165 "YYYY-MM-DD HH:MM:SS", style index 2
# in the original workbook) by copying the format from the row below, which
# already holds that style after the insert shifted it down.
$ws.Cells.Item(258, 4).Copy()
$ws.Range($ws.Cells.Item(256, 4), $ws.Cells.Item(257, 4)).PasteSpecial(-4122) # xlPasteFormats
$ws.Cells.Item(256, 4).Value = 44932
$ws.Cells.Item(257, 4).Value = 44932
$excel.CutCopyMode = 0
